# Applies the edits described by the commit: new test input values in
# column C (rows 8-25), the updated selected cell, and the minimized
# window state flag.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Window is shown minimized in the saved view.
$excel.Windows.Item(1).WindowState = -4140

# Update the raw "load" input values in column C; all dependent formulas
# (D:I, row 27, row 29, and the summary rows 31-36) recalculate automatically.
$ws.Range("C8").Value = 19
$ws.Range("C9").Value = 27
$ws.Range("C10").Value = 32
$ws.Range("C11").Value = 37
$ws.Range("C12").Value = 42
$ws.Range("C13").Value = 51
$ws.Range("C14").Value = 56
$ws.Range("C15").Value = 62
$ws.Range("C16").Value = 67
$ws.Range("C17").Value = 75
$ws.Range("C18").Value = 83
$ws.Range("C19").Value = 88
$ws.Range("C20").Value = 93
$ws.Range("C21").Value = 98
$ws.Range("C22").Value = 107
$ws.Range("C23").Value = 112
$ws.Range("C24").Value = 118
$ws.Range("C25").Value = 123

# Update the selected cell shown when the workbook is reopened.
$ws.Range("C32").Select()
